# "Generate Report for Handback"
#
# The localization-status report is regenerated after a handback event:
#  - The status text moves from "Ready for handoff" to
#    "Handed back: in sync with en-US" everywhere it appears (Overview +
#    both language sheets).
#  - Each language sheet (zh-cn, de-de) grows two new populated cells per
#    data row: "Latest Target File" (F) and "Latest Handback File" (G),
#    mirroring the existing "Source File Name" (A) / "Latest Handoff
#    File" (D) hyperlinked cells, since the handed-back file is the same
#    source file / xlf pair.
#  - "Latest Handback DateTime" (H) is stamped with the real handback
#    time instead of the zero-date placeholder. zh-cn and de-de were
#    handed back at different times, so they end up with different
#    strings there.

$wb = $excel.ActiveWorkbook

function Get-HyperlinkAddress($ws, $cellRef) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address(0, 0) -eq $cellRef) {
            return $hl.Address
        }
    }
    return $null
}

function Add-MirrorHyperlink($ws, $targetCellRef, $sourceCellRef, $displayText) {
    $addr = Get-HyperlinkAddress $ws $sourceCellRef
    $ws.Range($targetCellRef).Value = $displayText
    $ws.Hyperlinks.Add($ws.Range($targetCellRef), $addr, [Type]::Missing, [Type]::Missing, $displayText) | Out-Null
}

$handedBackStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status text for both rows/columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $handedBackStatus
$overview.Range("C2").Value = $handedBackStatus
$overview.Range("B3").Value = $handedBackStatus
$overview.Range("C3").Value = $handedBackStatus

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $handedBackStatus
$zhcn.Range("C3").Value = $handedBackStatus

# F = Latest Target File: both rows point at "a.md" (same source/address
#     as A2), matching D2/D3 which already both point at the same xlf
#     regardless of row.
# G = Latest Handback File (mirrors D, the handoff xlf + its hyperlink)
Add-MirrorHyperlink $zhcn "F2" "A2" "a.md"
Add-MirrorHyperlink $zhcn "G2" "D2" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
Add-MirrorHyperlink $zhcn "F3" "A2" "a.md"
Add-MirrorHyperlink $zhcn "G3" "D3" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# H = Latest Handback DateTime, now a real timestamp
$zhcn.Range("H2").Value = "2016-03-23 22:33:06"
$zhcn.Range("H3").Value = "2016-03-23 22:33:06"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $handedBackStatus
$dede.Range("C3").Value = $handedBackStatus

Add-MirrorHyperlink $dede "F2" "A2" "a.md"
Add-MirrorHyperlink $dede "G2" "D2" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
Add-MirrorHyperlink $dede "F3" "A2" "a.md"
Add-MirrorHyperlink $dede "G3" "D3" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$dede.Range("H2").Value = "2016-03-23 22:33:15"
$dede.Range("H3").Value = "2016-03-23 22:33:15"
